# Edit: insert 3 new data rows at row 274 of the single worksheet,
# shifting the existing data (rows 274:379) down to (277:382), and
# populate the new rows with a new "Frutilla" price-report group dated
# 2022-10-11 (serial 44845) for "Provincia de Melipilla".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 274. This pushes the
# existing rows 274:379 down to 277:382 and duplicates formatting
# (including the date number-format on column D) from the row above,
# matching how the original workbook is structured.
$ws.Range("A274:A276").EntireRow.Insert()

# Common / shared field values for the new group of three rows.
$mercadoId = 11
$mercado = "Vega Monumental Concepción"
$region = "Bíobío"
$fecha = 44845
$codreg = 8
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "`$/bandeja 7 kilos"
$origen = "Provincia de Melipilla"
$kgUnidad = 7

# Row-specific values: Calidad, Volumen, Precio min, Precio max, Precio prom, Precio $/Kg
$rows = @(
    @{ Row = 274; Calidad = "Especial"; Volumen = 100; PMin = 15000; PMax = 15000; PProm = 15000; PKg = 2143 },
    @{ Row = 275; Calidad = "Primera";  Volumen = 100; PMin = 13000; PMax = 13000; PProm = 13000; PKg = 1857 },
    @{ Row = 276; Calidad = "Segunda";  Volumen = 50;  PMin = 8000;  PMax = 8000;  PProm = 8000;  PKg = 1143 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
